# 23 Nov 2021 4th commit
# The "CCPAUrls" sheet tracks per-site CCPA opt-out toggles in column C
# ("RUN"). Rows 10-92 were still "OFF"; flip every one of them "ON".
# (Once nothing references the "OFF" string anymore, Excel drops it from
# the shared-string table on save, which is why uniqueCount shrinks too.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCPAUrls")
$ws.Activate() | Out-Null

$ws.Range("C10:C92").Value = "ON"

# Reproduce the view/selection state left behind in the saved sheet.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 71
$ws.Range("H83").Select() | Out-Null
